# Fix connector interface test fails
# Apply changes to the "SearchModelDataByCondition" sheet (sheet2):
#  - insert a new column F ("statusCode") shifting old F (expectCode) -> G
#    and old G (expectMessage) -> H
#  - set new F9 status code (400)
#  - update the "name does not exist" / "strange characters" / "too long"
#    rows (11-13) to use the new status code (108001) and new message text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchModelDataByCondition")

# Insert a new column before F, shifting expectCode/expectMessage right by one
$ws.Columns("F").Insert()

# New header for the inserted column
$ws.Range("F1").Value = "statusCode"

# Row 9 ("bad request (name is empty)") gets a new status code value
$ws.Range("F9").Value = 400
$ws.Range("H9").Value = "searchData.name is not valid,reason: must not be blank"

# Row 11 ("bad request (name does not exist)")
$ws.Range("G11").Value = 108001
$ws.Range("H11").Value = "The m2 service unavailable: (request M2 failed : no found entity )."

# Row 12 ("bad request (name contains strange characters)")
$ws.Range("G12").Value = 108001
$ws.Range("H12").Value = "The m2 service unavailable: (request M2 failed : no found entity )."

# Row 13 ("bad request (name is too long)")
$ws.Range("G13").Value = 108001
$ws.Range("H13").Value = "The m2 service unavailable: (request M2 failed : no found entity )."

# Update selection to match the authored workbook state
$ws.Range("H14").Select()
